# Auto-derived edit script: adds a "generator" component column to
# design_users and economics_users sheets, updates computed results.
$wb = $excel.ActiveWorkbook

# --- info_solution: updated computation time ---
$wsInfo = $wb.Worksheets.Item("info_solution")
$wsInfo.Range("A2").Value = 2.055131196975708

# --- design_users: insert "x_us_generator" column before x_us_batt (col F) ---
$wsDesign = $wb.Worksheets.Item("design_users")
$wsDesign.Range("F1").EntireColumn.Insert()

$wsDesign.Range("A1").Value = "User"
$wsDesign.Range("B1").Value = "Peak demand [kW]"
$wsDesign.Range("C1").Value = "Yearly Demand [MWh]"
$wsDesign.Range("D1").Value = "x_us_load"
$wsDesign.Range("E1").Value = "x_us_PV"
$wsDesign.Range("F1").Value = "x_us_generator"
$wsDesign.Range("G1").Value = "x_us_batt"
$wsDesign.Range("H1").Value = "x_us_conv"
$wsDesign.Range("I1").Value = "x_us_wind"
$wsDesign.Range("A2").Value = "user1"
$wsDesign.Range("B2").Value = 55.29324587
$wsDesign.Range("C2").Value = 232.10403242918818
$wsDesign.Range("E2").Value = 50.0
$wsDesign.Range("F2").Value = 0.0
$wsDesign.Range("A3").Value = "user2"
$wsDesign.Range("B3").Value = 27.00540954
$wsDesign.Range("C3").Value = 92.46745411323812
$wsDesign.Range("E3").Value = 25.489185845443384
$wsDesign.Range("G3").Value = 1.973245920368238
$wsDesign.Range("H3").Value = 1.973245920368238
$wsDesign.Range("A4").Value = "user3"
$wsDesign.Range("B4").Value = 45.40292054
$wsDesign.Range("C4").Value = 186.3047550336197
$wsDesign.Range("E4").Value = 44.04652629600476
$wsDesign.Range("G4").Value = 3.721368455391854
$wsDesign.Range("H4").Value = 3.721368455391854
$wsDesign.Range("I4").Value = 34.906983718279776

# --- economics_users: insert "C_gen_tot_us" (col E), "CAPEX_us_generator" (col M),
#     and "C_OEM_us_generator" (col S) columns ---
$wsEcon = $wb.Worksheets.Item("economics_users")
$wsEcon.Range("E1").EntireColumn.Insert()
$wsEcon.Range("M1").EntireColumn.Insert()
$wsEcon.Range("S1").EntireColumn.Insert()

$wsEcon.Range("A1").Value = "User_id"
$wsEcon.Range("B1").Value = "NPV_us"
$wsEcon.Range("C1").Value = "CAPEX_tot_us"
$wsEcon.Range("D1").Value = "yearly_rev"
$wsEcon.Range("E1").Value = "C_gen_tot_us"
$wsEcon.Range("F1").Value = "SDCF C_OEM_tot_us"
$wsEcon.Range("G1").Value = "SDCF C_REP_tot_us"
$wsEcon.Range("H1").Value = "SDCF R_RV_tot_us"
$wsEcon.Range("I1").Value = "SDCF C_Peak_tot_us"
$wsEcon.Range("J1").Value = "SDCF R_Energy_tot_us"
$wsEcon.Range("K1").Value = "CAPEX_us_load"
$wsEcon.Range("L1").Value = "CAPEX_us_PV"
$wsEcon.Range("M1").Value = "CAPEX_us_generator"
$wsEcon.Range("N1").Value = "CAPEX_us_batt"
$wsEcon.Range("O1").Value = "CAPEX_us_conv"
$wsEcon.Range("P1").Value = "CAPEX_us_wind"
$wsEcon.Range("Q1").Value = "C_OEM_us_load"
$wsEcon.Range("R1").Value = "C_OEM_us_PV"
$wsEcon.Range("S1").Value = "C_OEM_us_generator"
$wsEcon.Range("T1").Value = "C_OEM_us_batt"
$wsEcon.Range("U1").Value = "C_OEM_us_conv"
$wsEcon.Range("V1").Value = "C_OEM_us_wind"
$wsEcon.Range("A2").Value = "user1"
$wsEcon.Range("B2").Value = -563669.0832118867
$wsEcon.Range("C2").Value = 85000.0
$wsEcon.Range("D2").Value = -31580.287385073458
$wsEcon.Range("E2").Value = -0.00000000000008573955168024072
$wsEcon.Range("F2").Value = 22316.212290683256
$wsEcon.Range("G2").Value = 0.0
$wsEcon.Range("H2").Value = 9412.487821167684
$wsEcon.Range("I2").Value = 18246.63937566448
$wsEcon.Range("J2").Value = -447518.7193667072
$wsEcon.Range("L2").Value = 85000.0
$wsEcon.Range("M2").Value = 0.0
$wsEcon.Range("R2").Value = 1500.0
$wsEcon.Range("S2").Value = 0.0
$wsEcon.Range("A3").Value = "user2"
$wsEcon.Range("B3").Value = -223241.12057109879
$wsEcon.Range("C3").Value = 36868.807735841685
$wsEcon.Range("D3").Value = -12164.328747674768
$wsEcon.Range("E3").Value = 0.0
$wsEcon.Range("F3").Value = 11581.94006488829
$wsEcon.Range("G3").Value = 800.2766440442327
$wsEcon.Range("H3").Value = 4242.911954548019
$wsEcon.Range("I3").Value = 8840.453007913167
$wsEcon.Range("J3").Value = -169392.55507295928
$wsEcon.Range("L3").Value = 35684.86018362074
$wsEcon.Range("N3").Value = 789.2983681472952
$wsEcon.Range("O3").Value = 394.6491840736476
$wsEcon.Range("R3").Value = 764.6755753633015
$wsEcon.Range("T3").Value = 9.86622960184119
$wsEcon.Range("U3").Value = 3.946491840736476
$wsEcon.Range("A4").Value = "user3"
$wsEcon.Range("B4").Value = -435831.69036038255
$wsEcon.Range("C4").Value = 177428.21430168205
$wsEcon.Range("D4").Value = -17171.222950664833
$wsEcon.Range("E4").Value = 0.0
$wsEcon.Range("F4").Value = 35626.41777095829
$wsEcon.Range("G4").Value = 1509.2514460525526
$wsEcon.Range("H4").Value = 8353.446369498317
$wsEcon.Range("I4").Value = 9783.233210352739
$wsEcon.Range("J4").Value = -219838.02000083437
$wsEcon.Range("L4").Value = 70474.44207360761
$wsEcon.Range("N4").Value = 1488.5473821567416
$wsEcon.Range("O4").Value = 744.2736910783708
$wsEcon.Range("P4").Value = 104720.95115483932
$wsEcon.Range("R4").Value = 1321.3957888801426
$wsEcon.Range("T4").Value = 18.60684227695927
$wsEcon.Range("U4").Value = 7.442736910783708
$wsEcon.Range("V4").Value = 1047.2095115483933

# --- peak_users: recomputed rounding of Peak_id L value ---
$wsPeak = $wb.Worksheets.Item("peak_users")
$wsPeak.Range("O2").Value = 32.8611392682314

